# Updated cryptos list with latest Price / Volume(1h) figures, plus a
# ranking swap among ImmutableX / Filecoin / VeChain (rows 46-48).
#
# Price values in column D are free-text (e.g. "96.350.41", "0.0000251")
# rather than real numbers, so NumberFormat is forced to Text ("@") before
# assignment wherever the new value would otherwise be auto-coerced to a
# number by Excel (which would also mangle things like trailing zeros and
# very small magnitudes into scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.350.41'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '3.331.91'
$ws.Range("E3").Value = '  -2.88%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.63'
$ws.Range("E5").Value = '  -2.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '654.25'
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.38'
$ws.Range("E7").Value = '  -7.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.421'
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.996'
$ws.Range("E10").Value = '  -6.15%  '
$ws.Range("D11").Value = '3.330.12'
$ws.Range("E11").Value = '  -2.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.206'
$ws.Range("E12").Value = '  -3.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.25'
$ws.Range("E13").Value = '  -4.36%  '
$ws.Range("D14").Value = '96.095.35'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.08'
$ws.Range("E15").Value = '  -5.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000251'
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("D17").Value = '3.952.38'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.49'
$ws.Range("E18").Value = '  -2.94%  '
$ws.Range("D19").Value = '3.328.53'
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.04'
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.508'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '503.14'
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.53'
$ws.Range("E23").Value = '  -5.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.34'
$ws.Range("E24").Value = '  -3.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000198'
$ws.Range("E25").Value = '  -4.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.54'
$ws.Range("E26").Value = '  +7.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.98'
$ws.Range("E27").Value = '  -2.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.07'
$ws.Range("E28").Value = '  -6.06%  '
$ws.Range("D29").Value = '3.502.51'
$ws.Range("E29").Value = '  -2.86%  '
$ws.Range("E30").Value = '  -7.24%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.08'
$ws.Range("E32").Value = '  -3.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.186'
$ws.Range("E33").Value = '  -6.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.47'
$ws.Range("E34").Value = '  +8.99%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.545'
$ws.Range("E36").Value = '  -5.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.01'
$ws.Range("E37").Value = '  -7.00%  '
$ws.Range("E38").Value = '  +3.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.60'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -3.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '508.91'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.833'
$ws.Range("E44").Value = '  -4.46%  '
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0413'
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.67'
$ws.Range("E47").Value = '  +5.53%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.49'
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.35'
$ws.Range("E49").Value = '  +1.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.11'
$ws.Range("E50").Value = '  +3.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.13'
$ws.Range("E51").Value = '  -5.53%  '
